# Update API doc sheet (Sheet1) to reflect the new "concat" behaviour for
# text resources and simplify the "function" resource description, and move
# the active selection to B8 (matches the author's last edit position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "concat original and posted strings"
$ws.Range("B6").Value = "return function description"

$ws.Range("B8").Select()
